# Normalize the "Recorded By" (column G) entries so that within each
# comma-separated list of recorders, the first two names are swapped -
# unless they already read "System, backup@backdoor.com" (already in the
# canonical order). This mirrors the upstream sync that reordered the
# recorder names in attendance_reports.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$lastCol = $ws.UsedRange.Columns.Count

# Column G = "Recorded By" (7th column); confirm via header just in case
$recordedByCol = 7
$header = $ws.Cells.Item(1, $recordedByCol).Value2
if ($header -ne "Recorded By") {
    for ($c = 1; $c -le $lastCol; $c++) {
        if ($ws.Cells.Item(1, $c).Value2 -eq "Recorded By") {
            $recordedByCol = $c
        }
    }
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        if ($parts.Count -ge 2) {
            $alreadyCanonical = ($parts[0] -eq "System" -and $parts[1] -eq "backup@backdoor.com")

            if (-not $alreadyCanonical) {
                $first = $parts[0]
                $second = $parts[1]
                $parts[0] = $second
                $parts[1] = $first

                $newVal = [string]::Join(", ", $parts)
                $cell.Value2 = $newVal
            }
        }
    }
}
